# Insert a new weekly price observation row for
# "Femacal de La Calera - Haba" above the current row 276,
# shifting the existing rows 276-285 down to 277-286.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 276 (pushes rows 276..285 down to 277..286,
# carrying formatting/styles from the row being pushed down).
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(276, 1).Value = 3
$ws.Cells.Item(276, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(276, 3).Value = "Coquimbo"
$ws.Cells.Item(276, 4).Value = 45147
$ws.Cells.Item(276, 5).Value = 5
$ws.Cells.Item(276, 6).Value = 100112026
$ws.Cells.Item(276, 7).Value = "Haba"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 40
$ws.Cells.Item(276, 11).Value = 15000
$ws.Cells.Item(276, 12).Value = 15000
$ws.Cells.Item(276, 13).Value = 15000
$ws.Cells.Item(276, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(276, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(276, 16).Value = 600
$ws.Cells.Item(276, 17).Value = 25
$ws.Cells.Item(276, 18).Value = "Hortaliza"
